$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:K1")
$rng.NumberFormat = "@"

$ws.Range("A1").Value = "Rush"
$ws.Range("B1").Value = "Cooper"
$ws.Range("C1").Value = "QB"
$ws.Range("D1").Value = "2018-10-14"
$ws.Range("E1").Value = "6"
$ws.Range("F1").Value = "24.327"
$ws.Range("G1").Value = "DAL"
$ws.Range("H1").Value = ""
$ws.Range("I1").Value = "JAX"
$ws.Range("J1").Value = "W 40-7"
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = 0
